$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.235.23'
$ws.Range("E2").Value = '  -3.00%  '
$ws.Range("D3").Value = '3.296.03'
$ws.Range("E3").Value = '  -3.69%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -3.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '140.65'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -8.72%  '
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("D8").Value = '3.296.39'
$ws.Range("E8").Value = '  -3.66%  '
$ws.Range("E10").Value = '  -2.48%  '
$ws.Range("E11").Value = '  -5.33%  '
$ws.Range("E12").Value = '  -2.80%  '
$ws.Range("D13").Value = '3.861.67'
$ws.Range("E13").Value = '  -3.66%  '
$ws.Range("E14").Value = '  -0.32%  '
$ws.Range("D16").Value = '3.295.52'
$ws.Range("E16").Value = '  -3.72%  '
$ws.Range("E17").Value = '  -5.00%  '
$ws.Range("D18").Value = '60.220.91'
$ws.Range("E18").Value = '  -3.05%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.06'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -7.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.61'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.69%  '
$ws.Range("E21").Value = '  -5.20%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '373.10'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.71%  '
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '71.94'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -5.26%  '
$ws.Range("E25").Value = '  -7.13%  '
$ws.Range("D26").Value = '3.434.70'
$ws.Range("E26").Value = '  -3.62%  '
$ws.Range("E27").Value = '  -9.63%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.173'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.29%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("E30").Value = '  -8.09%  '
$ws.Range("E31").Value = '  -0.09%  '
$ws.Range("E32").Value = '  -5.01%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '7.32'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -6.94%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.56'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.06%  '
$ws.Range("E35").Value = '  -4.98%  '
$ws.Range("B36").Value = 'Monero'
$ws.Range("C36").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '165.59'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.87%  '
$ws.Range("B37").Value = 'NEARProtocol'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.02'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -9.10%  '
$ws.Range("E38").Value = '  -4.46%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.62'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -4.90%  '
$ws.Range("D40").Value = '3.328.41'
$ws.Range("E40").Value = '  -3.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0720'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '25.50'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -17.29%  '
$ws.Range("E43").Value = '  -2.41%  '
$ws.Range("E44").Value = '  -4.53%  '
$ws.Range("E45").Value = '  -4.23%  '
$ws.Range("E46").Value = '  -7.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.73%  '
$ws.Range("E48").Value = '  -0.04%  '
$ws.Range("D49").Value = '2.321.39'
$ws.Range("E49").Value = '  -9.20%  '
$ws.Range("E50").Value = '  -6.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '21.43'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -7.92%  '
